$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 12 (20020606/274.0), shifting
# everything below down by two rows.
$ws.Rows("12:13").Insert()

# Force text formatting on the two new date cells so they are stored the
# same way as every other date in column A (text, not a number), then
# write the values.
$ws.Range("A12:A13").NumberFormat = "@"

$ws.Range("A12").Value = "20020403"
$ws.Range("B12").Value = 310

$ws.Range("A13").Value = "20020411"
$ws.Range("B13").Value = 343
